$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New work log rows 12-14. Row 8 already carries the exact style pattern we
# need (A=s2 NAME style, C=s4 date style, E=s2 plain, G=s5 hyperlink style,
# H=s2 plain text), so copy formats-only from row 8 into the new rows first.
# ---------------------------------------------------------------------------
foreach ($row in 12, 13, 14) {
    foreach ($col in "A", "C", "E", "G", "H") {
        $ws.Range("$col`8").Copy() | Out-Null
        $ws.Range("$col$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    }
}

# ---------------------------------------------------------------------------
# Enter the cell text/values. Order matches how the shared-string table was
# actually built up (row 13's commit/note first, then the name, then row 12,
# then row 14).
# ---------------------------------------------------------------------------

# Row 13
$ws.Range("G13").Value = "https://github.com/chat-loc/chatloc.github.io/commit/a28fef4d8582293946dbff07fa83c1827432158f"
$ws.Range("H13").Value = "Creating chat user interface"
$ws.Range("A13").Value = "Khushboo"
$ws.Range("C13").Value = 44032
$ws.Range("E13").Value = 400

# Row 12
$ws.Range("G12").Value = "https://github.com/chat-loc/chatloc.github.io/pull/26/commits/3398a9cf0b0f7920890fe363bfc6e728045f48c7"
$ws.Range("H12").Value = "Created landing page, css files and login page"
$ws.Range("A12").Value = "Khushboo"
$ws.Range("C12").Value = 44023
$ws.Range("E12").Value = 450

# Row 14
$ws.Range("G14").Value = "https://github.com/chat-loc/chatloc.github.io/commit/a45398dcb5d9e7db44ce3e49223b552742ae46c3"
$ws.Range("H14").Value = "Making login logout and room list responsive"
$ws.Range("A14").Value = "Khushboo"
$ws.Range("C14").Value = 44014
$ws.Range("E14").Value = 390

# ---------------------------------------------------------------------------
# Wire up the G-column hyperlinks (text already matches the target URL, so
# no TextToDisplay override is needed -- that would otherwise clobber the
# cell text). Added in G12, G13, G14 order so the relationship ids line up
# with the document order of the <hyperlinks> entries.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G12"), $ws.Range("G12").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G13"), $ws.Range("G13").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G14"), $ws.Range("G14").Value2) | Out-Null

# Hyperlinks.Add() resets the cell style, so reapply the row-8 G-column style.
foreach ($row in 12, 13, 14) {
    $ws.Range("G8").Copy() | Out-Null
    $ws.Range("G$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Leave the selection where the author left off, just past the new rows.
$ws.Range("H16").Select() | Out-Null
